$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("华胜天成", "浙江世宝", "浙江世宝"),
    @("天奇股份", "华胜天成", "华胜天成"),
    @("浙江世宝", "深科技", "汉缆股份"),
    @("长飞光纤", "紫金矿业", "天奇股份"),
    @("掌阅科技", "天奇股份", "利欧股份"),
    @("风语筑", "华工科技", "风语筑"),
    @("紫金矿业", "三花智控", "博纳影业"),
    @("深科技", "长飞光纤", "三花智控"),
    @("三花智控", "风语筑", "巨力索具"),
    @("汉缆股份", "东方财富", "紫金矿业"),
    @("利欧股份", "贵州茅台", "长飞光纤"),
    @("光线传媒", "光线传媒", "掌阅科技"),
    @("电广传媒", "掌阅科技", "大位科技"),
    @("华工科技", "利欧股份", "深科技"),
    @("百达精工", "汉缆股份", "光线传媒"),
    @("万向钱潮", "百达精工", "万向钱潮"),
    @("东方财富", "电广传媒", "嘉美包装"),
    @("昆仑万维", "赣锋锂业", "二六三"),
    @("贵州茅台", "万向钱潮", "百达精工"),
    @("寒武纪-U", "亚星锚链", "克来机电")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
